$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the encumbrance-type code used in row 2 (pledge -> pledge in favor of the bank itself)
$ws.Range("C2").Value = "PLEDGE_OUR"

# Reflect the active-cell selection recorded in the saved workbook
$ws.Range("C2").Select()
